$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "last updated" timestamp title
$ws.Range("A1").Value = "Datos actualizados a 8 de Mayo de 2020 a las 23:34"

# Update Estados Unidos (row 4) totals
$ws.Cells.Item(4, 2).Value = 1318504
$ws.Cells.Item(4, 3).Value = 25881
$ws.Cells.Item(4, 4).Value = 221919
$ws.Cells.Item(4, 5).Value = 1018093
$ws.Cells.Item(4, 6).Value = 16928
$ws.Cells.Item(4, 7).Value = 1564
$ws.Cells.Item(4, 8).Value = 78492

# Other independent numeric updates
$ws.Cells.Item(20, 6).Value = 141   # row20 F (Peru)
$ws.Cells.Item(22, 4).Value = 26100 # row22 D
$ws.Cells.Item(22, 5).Value = 2284  # row22 E
$ws.Cells.Item(22, 6).Value = 101   # row22 F
$ws.Cells.Item(22, 7).Value = 13    # row22 G
$ws.Cells.Item(22, 8).Value = 1823  # row22 H
$ws.Cells.Item(26, 6).Value = 508   # row26 F
$ws.Cells.Item(107, 4).Value = 537  # row107 D
$ws.Cells.Item(107, 5).Value = 168  # row107 E

# Reorder "Sudan del Sur" to just after Camboya (row 155), shifting Bermudas..Bahamas down by one, with refreshed stats
$ws.Cells.Item(156, 1).Value = "Sudan del Sur"
$ws.Cells.Item(156, 2).Value = 120
$ws.Cells.Item(156, 3).Value = 46
$ws.Cells.Item(156, 4).Value = 2
$ws.Cells.Item(156, 5).Value = 118
$ws.Cells.Item(156, 6).Value = 0
$ws.Cells.Item(156, 7).Value = 0
$ws.Cells.Item(156, 8).Value = 0

$ws.Cells.Item(157, 1).Value = "Bermudas"
$ws.Cells.Item(157, 2).Value = 118
$ws.Cells.Item(157, 3).Value = 0
$ws.Cells.Item(157, 4).Value = 61
$ws.Cells.Item(157, 5).Value = 50
$ws.Cells.Item(157, 6).Value = 4
$ws.Cells.Item(157, 7).Value = 0
$ws.Cells.Item(157, 8).Value = 7

$ws.Cells.Item(158, 1).Value = "Trinidad yTobago"
$ws.Cells.Item(158, 2).Value = 116
$ws.Cells.Item(158, 3).Value = 0
$ws.Cells.Item(158, 4).Value = 103
$ws.Cells.Item(158, 5).Value = 5
$ws.Cells.Item(158, 6).Value = 0
$ws.Cells.Item(158, 7).Value = 0
$ws.Cells.Item(158, 8).Value = 8

$ws.Cells.Item(159, 1).Value = "Nepal"
$ws.Cells.Item(159, 2).Value = 102
$ws.Cells.Item(159, 3).Value = 1
$ws.Cells.Item(159, 4).Value = 31
$ws.Cells.Item(159, 5).Value = 71
$ws.Cells.Item(159, 6).Value = 0
$ws.Cells.Item(159, 7).Value = 0
$ws.Cells.Item(159, 8).Value = 0

$ws.Cells.Item(160, 1).Value = "Uganda"
$ws.Cells.Item(160, 2).Value = 101
$ws.Cells.Item(160, 3).Value = 0
$ws.Cells.Item(160, 4).Value = 55
$ws.Cells.Item(160, 5).Value = 46
$ws.Cells.Item(160, 6).Value = 0
$ws.Cells.Item(160, 7).Value = 0
$ws.Cells.Item(160, 8).Value = 0

$ws.Cells.Item(161, 1).Value = "Aruba"
$ws.Cells.Item(161, 2).Value = 101
$ws.Cells.Item(161, 3).Value = 0
$ws.Cells.Item(161, 4).Value = 89
$ws.Cells.Item(161, 5).Value = 9
$ws.Cells.Item(161, 6).Value = 4
$ws.Cells.Item(161, 7).Value = 0
$ws.Cells.Item(161, 8).Value = 3

$ws.Cells.Item(162, 1).Value = "Monaco"
$ws.Cells.Item(162, 2).Value = 95
$ws.Cells.Item(162, 3).Value = 0
$ws.Cells.Item(162, 4).Value = 82
$ws.Cells.Item(162, 5).Value = 9
$ws.Cells.Item(162, 6).Value = 1
$ws.Cells.Item(162, 7).Value = 0
$ws.Cells.Item(162, 8).Value = 4

$ws.Cells.Item(163, 1).Value = "Republica de Africa Central"
$ws.Cells.Item(163, 2).Value = 94
$ws.Cells.Item(163, 3).Value = 0
$ws.Cells.Item(163, 4).Value = 10
$ws.Cells.Item(163, 5).Value = 84
$ws.Cells.Item(163, 6).Value = 0
$ws.Cells.Item(163, 7).Value = 0
$ws.Cells.Item(163, 8).Value = 0

$ws.Cells.Item(164, 1).Value = "Guyana"
$ws.Cells.Item(164, 2).Value = 94
$ws.Cells.Item(164, 3).Value = 1
$ws.Cells.Item(164, 4).Value = 34
$ws.Cells.Item(164, 5).Value = 50
$ws.Cells.Item(164, 6).Value = 3
$ws.Cells.Item(164, 7).Value = 0
$ws.Cells.Item(164, 8).Value = 10

$ws.Cells.Item(165, 1).Value = "Bahamas"
$ws.Cells.Item(165, 2).Value = 92
$ws.Cells.Item(165, 3).Value = 0
$ws.Cells.Item(165, 4).Value = 31
$ws.Cells.Item(165, 5).Value = 50
$ws.Cells.Item(165, 6).Value = 1
$ws.Cells.Item(165, 7).Value = 0
$ws.Cells.Item(165, 8).Value = 11

# Reorder "Angola" to just after Macao (row 174), shifting Malaui..San Martin (Parte Francesa) down by one, with refreshed stats
$ws.Cells.Item(175, 1).Value = "Angola"
$ws.Cells.Item(175, 2).Value = 43
$ws.Cells.Item(175, 3).Value = 7
$ws.Cells.Item(175, 4).Value = 11
$ws.Cells.Item(175, 5).Value = 30
$ws.Cells.Item(175, 6).Value = 0
$ws.Cells.Item(175, 7).Value = 0
$ws.Cells.Item(175, 8).Value = 2

$ws.Cells.Item(176, 1).Value = "Malaui"
$ws.Cells.Item(176, 2).Value = 43
$ws.Cells.Item(176, 3).Value = 0
$ws.Cells.Item(176, 4).Value = 14
$ws.Cells.Item(176, 5).Value = 26
$ws.Cells.Item(176, 6).Value = 1
$ws.Cells.Item(176, 7).Value = 0
$ws.Cells.Item(176, 8).Value = 3

$ws.Cells.Item(177, 1).Value = "Mongolia"
$ws.Cells.Item(177, 2).Value = 42
$ws.Cells.Item(177, 3).Value = 1
$ws.Cells.Item(177, 4).Value = 13
$ws.Cells.Item(177, 5).Value = 29
$ws.Cells.Item(177, 6).Value = 0
$ws.Cells.Item(177, 7).Value = 0
$ws.Cells.Item(177, 8).Value = 0

$ws.Cells.Item(178, 1).Value = "Puerto Rico"
$ws.Cells.Item(178, 2).Value = 39
$ws.Cells.Item(178, 3).Value = 0
$ws.Cells.Item(178, 4).Value = 1
$ws.Cells.Item(178, 5).Value = 36
$ws.Cells.Item(178, 6).Value = 0
$ws.Cells.Item(178, 7).Value = 0
$ws.Cells.Item(178, 8).Value = 2

$ws.Cells.Item(179, 1).Value = "Eritrea"
$ws.Cells.Item(179, 2).Value = 39
$ws.Cells.Item(179, 3).Value = 0
$ws.Cells.Item(179, 4).Value = 37
$ws.Cells.Item(179, 5).Value = 2
$ws.Cells.Item(179, 6).Value = 0
$ws.Cells.Item(179, 7).Value = 0
$ws.Cells.Item(179, 8).Value = 0

$ws.Cells.Item(180, 1).Value = "San Martin (Parte Francesa)"
$ws.Cells.Item(180, 2).Value = 38
$ws.Cells.Item(180, 3).Value = 0
$ws.Cells.Item(180, 4).Value = 30
$ws.Cells.Item(180, 5).Value = 5
$ws.Cells.Item(180, 6).Value = 1
$ws.Cells.Item(180, 7).Value = 0
$ws.Cells.Item(180, 8).Value = 3
